$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "2019.03.09"
$ws.Range("C15").Value = "#jupyter"
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = "Chapter one math"

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "2019.03.09"
$ws.Range("C16").Value = "#jupyter"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = "Speed triangles"

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "2019.03.09"
$ws.Range("C17").Value = "#latex"
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = "Speed triangles"

$ws.Range("D21").Select()
